{"js": "// Fix the Exercise 1 grading percentage: \"25% ex 1\" -> \"20% ex 1\"\n// (the rest of the paragraph, \"20% ex 2 respectively\", already reads 20%\n// and is left untouched).\nconst results = context.document.body.search(\"25% ex 1\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target text \"25% ex 1\" not found in document body.');\n}\n\n// Replace in place so formatting / run properties of the existing run are preserved.\nresults.items[0].insertText(\"20% ex 1\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix the Exercise 1 grading percentage: \"25% ex 1\" -> \"20% ex 1\"\n# (the rest of the paragraph, \"20% ex 2 respectively\", already reads 20%\n# and is left untouched).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"25% ex 1\"\n$find.Replacement.Text = \"20% ex 1\"\n$find.Execute(\n    \"25% ex 1\",   # FindText\n    $true,        # MatchCase\n    $false,       # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    \"20% ex 1\",   # ReplaceWith\n    2             # Replace (wdReplaceAll)\n)\n"}
